$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: renumber id column (A) for existing rows 2-85 to be sequential (row-1) ---
for ($r = 2; $r -le 85; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- Step 2: column width adjustments ---
$ws.Columns.Item(4).ColumnWidth = 28.136   # raid_id column, now holds long quest-group names
$ws.Columns.Item(5).ColumnWidth = 23.423   # required_quest_id column
$ws.Columns.Item(28).ColumnWidth = 28.136  # new AB column (secondary raid_id)
$ws.Columns.Item(29).ColumnWidth = 23.423  # new AC column (secondary required_quest_id)

# --- Step 3: add new quest rows 86-96 ("Pirates and Merchants" questline) ---

# Row 86
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = 'Pirates and Merchants'
$ws.Cells.Item(86, 3).Value = 'TheSoldier'
$ws.Cells.Item(86, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(86, 7).Value = 'Stolen Coin'
$ws.Cells.Item(86, 15).Value = 50000000
$ws.Cells.Item(86, 16).Value = 1000
$ws.Cells.Item(86, 19).Value = 1
$ws.Cells.Item(86, 26).Value = '“Listen child, have you heard the story of the Pirates and the Merchants?” The Soldier began. The two of you sat together in an INN not far from the docs.<br /> <br /> “There was a time this place, Smugglers Port was run by Pirates. There were merchants that roamed the land in caravans and even a few dared to venture into the walls of the then pirate run city of ramshackle huts and whore houses.”<br /> <br /> He paused to take a gander out the window the two of you sat in front of at a table, drinking ale and eating various meats.<br /> <br /> “It is said that one the merchants who came into the city, stole from the Pirate Lord who now roams the lands looking for the lost item – which then started a war. The specifics of what it was that was stolen alludes me, maybe it was a coin?”'
$ws.Cells.Item(86, 27).Value = 'You slaughtered the foes before you, and in the blood soaked ground, buried slightly under the upturned soil of a fight, lies a glimmering shiny object. You pull it out and clean it off, it’s a coin. Maybe the same one The Soldier was talking about.<br /> <br /> Before your eyes, a flash of light, a vision of sorts, a story unfolds, a drunkard stands in the alley pissing on a wall.<br /> <br /> You know the face, familiar.<br /> <br /> The man finishes and heads deeper into the alley. You follow down the long dark alley, deeper into the shadows. Another man.<br /> <br /> The two men embrace in a lovers passionate kiss in the shadows of the alley.<br /> <br /> “Child!” You open your eyes and see the Soldier Standing over you.<br /> <br /> “Whats that in your hand?”'
$ws.Cells.Item(86, 28).Value = 'The Smugglers Are Back!'

# Row 87
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = 'Trade Ya'
$ws.Cells.Item(87, 3).Value = 'ThePoet'
$ws.Cells.Item(87, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(87, 7).Value = 'Smugglers Coin'
$ws.Cells.Item(87, 12).Value = 'Pirates Lost Key'
$ws.Cells.Item(87, 15).Value = 50000000
$ws.Cells.Item(87, 16).Value = 1000
$ws.Cells.Item(87, 20).Value = 'Pirates and Merchants'
$ws.Cells.Item(87, 26).Value = '“You spent a week staring at this coin, wondering what it means. The dream or vision plays secondary on your mind. You play it out, the familiarity of the face, but you could not place it.”<br /> <br /> You go searching for The Poet, thinking he might know more. He is wise after all, had something to do with The Creator.<br /> <br /> “I am unsure I know the story you speak of child,” The Poet states as you managed to track him down at the local market and chat him up about the coin and what The Soldier told you.<br /> <br /> “I do know a war between the Pirates and the Merchants did take place and that’s why Smugglers Port is run by the High Merchant Council.”<br /> <br /> The two of you continue to walk through the market place talking more about the war. He studies the coin you found and states he saw one just like it a while ago that the Drunk Adventurer had in his possession and would break about how he fund a coin for his lover in the dungeons of despair.'
$ws.Cells.Item(87, 27).Value = 'Hacking, slashing, burning through the foes before you, they fall one after the other. You see another coin in the mud and pick it up. This one seems similar to the other one. Holding it gives you the feelings of love, undeniable, lustful and shameless provocative love.<br /> <br /> You can see the INN room, the bodies writhing to the rhythm of their own fast beating heart beats. You can see the shadows dance the dance of love under the night sky, kept alive only by a single burning flame.<br /> You snap back to reality. Why are you visiting these two people while they are acting like dogs in heat? What does it have to do with the coins?<br /> <br /> “I shall take that child.”<br /> <br /> Take the coin? Why? You are the one that found them. You want the other one back.<br /> <br /> “There is a magic in those coins. A curse. A vision of the past. One coin is the other and the other is the pirates coin.”<br /> Is it that simple? A merchant and a pirate? A war over love?'
$ws.Cells.Item(87, 28).Value = 'The Smugglers Are Back!'

# Row 88
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = 'Keys to the past'
$ws.Cells.Item(88, 3).Value = 'LabyrinthOracle'
$ws.Cells.Item(88, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(88, 7).Value = 'Pirates Lost Key'
$ws.Cells.Item(88, 12).Value = 'Corrupted Bracelet'
$ws.Cells.Item(88, 15).Value = 20000000
$ws.Cells.Item(88, 16).Value = 1000
$ws.Cells.Item(88, 20).Value = 'Trade Ya'
$ws.Cells.Item(88, 26).Value = 'You walk for days through the mazes of the Labyrinth plane. You see out the oracle for the answers to your questions. The Key, it wont stop showing you the images of death. Of war, of pain and suffering. You feel the projected emotional pain of the killer.<br /> <br /> Yet he turns around and fucks his lover with passion, with love, with vigor. Blood stains his body, his hands, his lips. Is his lover none the wiser?<br /> <br /> “I see you wondering around in the labyrinth of your own mind child.” A voice comes. The oracle stands before you.<br /> <br /> “You seek to rid your self of the cure of the truth are you?” You ask her, beg of her anything she can do. Explain this. Why?!<br /> <br /> “Hush child as I tell you of love, war and sacrifice. The short of it is simple: The Pirate Lord has a child, this child grew to love another child. A child of a powerful Merchant. Love is complicated, fickle and finicky. Love is raw and passionate and unending. Love has no master, no will to bend.”<br /> <br /> The Oracle pauses and takes the key from your hand.<br /> <br /> “Child, the voices are gone, but you have more answers to find.”<br /> <br /> It seems pretty obvious that this is nothing more then forbidden love and that wars have been technically fought for less.<br /> <br /> How ever what you find upon your arm now gives you a glimpse of tomorrow.'
$ws.Cells.Item(88, 27).Value = 'There is a scream from the tent, pain, agony, blood curdling scream.<br /> <br /> He races towards the tent, only to find his lover’s body broken and scattered. The blood covers the walls of the tent, all he can do is scream in terror.<br /> <br /> You see the camp burn, you see the men on horses race through, race by and stick their swords in anything that moved. There eyes are as blue as the day is bright. Their faces are familiar. They are not Red Hawks. Something about the Sky Judges.<br /> <br /> Your mind wonders to fall and you wind up back in reality with a man staring at you. The Merchant.'
$ws.Cells.Item(88, 28).Value = 'The Smugglers Are Back!'

# Row 89
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = 'Loves is an after thought'
$ws.Cells.Item(89, 3).Value = 'WanderingMerchant'
$ws.Cells.Item(89, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(89, 7).Value = 'Corrupted Bracelet'
$ws.Cells.Item(89, 8).Value = 45000
$ws.Cells.Item(89, 15).Value = 10000000
$ws.Cells.Item(89, 16).Value = 1500
$ws.Cells.Item(89, 20).Value = 'Keys to the past'
$ws.Cells.Item(89, 21).Value = 'Shade Dust'
$ws.Cells.Item(89, 26).Value = '“Hello child. Fancy seeing you again. I have a favor to ask you.” You stare at him for a moment. He is no longer old and haggard. He is young and handsome.<br /> <br /> You bluntly ask him if he is the pirates lover and in so doing you realize what you have asked and cover your mouth as quickly as possible, embarrassed and ashamed.<br /> <br /> “No Child I am not the lover of this Pirate you seek. Alas he is slaughtering people in Smuggler Port. It sort of makes sense, given the fact that his lover was murdered. Now do I look dead?”<br /> <br /> You tell him no and that he doesn’t look old and haggard either.<br /> <br /> “Cosmic beings beyond your pay grade rarely do.”<br /> <br /> You stare at him. So many answers.<br /> <br /> “Listen to me. I need you to get his lovers dust, also known as Shade Dust. If I mix it with the corrupted bracelet I can create a portal to the past, for a brief moment in time the answer might be made clear, but not in the way you might expect.”<br /> <br /> Great more riddles. Riddles do lead to answers though.'
$ws.Cells.Item(89, 27).Value = 'You have the dust of the shade known as the lover. The merchant. This story is confusing, there is love, there is death, there is no reason though. It seems like a senseless war over the death of a pirate lords son, which in it’s self is technically a reason to start a war, but both sides knew a free form expression of love was not tolerated under their then ruler: The Federation.<br /> <br /> The One True God, the Creator of all that there is forbid it. The Children were beaten for it, the lovers were persecuted, they were slaughtered before the people. The Church.<br /> You hand the dust and the bracelet over to The Merchant and he concocts a spell that opens a portal before your eyes.<br /> <br /> “A brief moment child. Step through now.”<br /> <br /> You walk into the portal, it glows, it closes and you can see the truth before you. The answers lie with in your eye sight. The death that rained down on the camp, the war that was started for love, cliche as it is, was never started by any one other then those seeking retribution.<br /> <br /> You stand before Rhe Merchant once more and he states you know what you have to do.'
$ws.Cells.Item(89, 28).Value = 'The Smugglers Are Back!'

# Row 90
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = 'It Was All A Lie'
$ws.Cells.Item(90, 3).Value = 'ScholarofHell'
$ws.Cells.Item(90, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(90, 7).Value = 'Scholars Necklace'
$ws.Cells.Item(90, 15).Value = 15000000
$ws.Cells.Item(90, 16).Value = 1500
$ws.Cells.Item(90, 20).Value = 'Loves is an after thought'
$ws.Cells.Item(90, 21).Value = 'The Quartermasters Severed Head'
$ws.Cells.Item(90, 26).Value = 'You went searching and you found him. Standing at the end of a lake in Hell. The Scholar.<br /> <br /> You beckoned to him and he turned and smiled with a wicked grin.<br /> <br /> You approach him and come closer, but he Is demonic in nature and slips in to the burning fires of the lake before him. His flesh melts off and his bones turn to dust.<br /> Before you stands a empty open chest.'
$ws.Cells.Item(90, 27).Value = 'You managed to figure it out. The Quarter master is the one who tipped them off, he is the one who said they would make a perfect target to start a war, one that would ultimately draw The Church in. The set of individual who come to the aid of The Merchants and push these savages out of town.<br /> <br /> The question that remained was why? Why turn on your own?<br /> <br /> “Because the strike was never meant for the merchant lover child. I thought you would have gathered as much by now.”<br /> <br /> Never meant for the merchant, it was meant for the son, it was designed as a trojan horse with out being one. But it back fired. The war still happened, it was still bloody, and they came, and helped. They slaughtered The Pirates.<br /> <br /> The Church did come, the soldiers of God wanted to burn the city to the ground, but when they showed up, it was a blood bath. They who went along with this plan.<br /> <br /> Why when The Creator was apart of their squad? Why when he had his own love story end in tragedy? Was this what corrupted The Creator, turned him into an entity unknown, although all powerful? A wrathful and angry child?'
$ws.Cells.Item(90, 28).Value = 'The Smugglers Are Back!'

# Row 91
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = 'The Creator and The Hawks'
$ws.Cells.Item(91, 3).Value = 'ScholarofHell'
$ws.Cells.Item(91, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(91, 7).Value = 'Simple Gem Necklace'
$ws.Cells.Item(91, 8).Value = 3000
$ws.Cells.Item(91, 12).Value = 'Merchants Walking Stick'
$ws.Cells.Item(91, 15).Value = 5000000
$ws.Cells.Item(91, 16).Value = 1600
$ws.Cells.Item(91, 20).Value = 'It Was All A Lie'
$ws.Cells.Item(91, 26).Value = 'You remembered a story once told told to you by the Scholar him self. A story of how The Creator escaped Purgatory but Ascended with rage and anger and became a manifestation of evil. Of pain, sorrow and hatred that filled his very being. It was sad really, if that was indeed the case, because of how at the end of the day this “being of immense power” is nothing more then a broken heart ed man with self inflected wounds and a love he cannot let go of.<br /> <br /> The Chest stands empty once again.'
$ws.Cells.Item(91, 27).Value = 'You place the items in the chest and stand before the burning lake watching and waiting for any sign of the Scholar.<br /> A voice that whispers to your thoughts haunts and taunts you with its soft and supple words.<br /> <br /> “The Red Hawks”<br /> <br /> An organization which fought against The Federation and the Church. An organization that liberated and freed the oppressed people of the North. They have freed lands never before seen before. The Child was found by them, guided by them under the instruction of The Poet.<br /> <br /> They initiated the attack, they had the wrong target, they had … the wrong target. That''s what started the war but the fact that they would kill one of the people they wanted to liberated from the oppression of The Church is what shakes you. They sacrificed people.<br /> <br /> “Oh grow up child” the Merchant stands before you.'
$ws.Cells.Item(91, 28).Value = 'The Smugglers Are Back!'

# Row 92
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = 'The Creators Choice'
$ws.Cells.Item(92, 3).Value = 'WanderingMerchant'
$ws.Cells.Item(92, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(92, 7).Value = 'Merchants Walking Stick'
$ws.Cells.Item(92, 8).Value = 25000
$ws.Cells.Item(92, 13).Value = 1000
$ws.Cells.Item(92, 14).Value = 100
$ws.Cells.Item(92, 15).Value = 25000000
$ws.Cells.Item(92, 16).Value = 3000
$ws.Cells.Item(92, 20).Value = 'The Creator and The Hawks'
$ws.Cells.Item(92, 26).Value = '“You need to realize the world is a much more complicated one then the fantasy you got cooking up in your head.”<br /> <br /> The Merchant pauses and walks over to a rock. The world of the surface is back in your view. The trees are lush and green and the sky is blue and clear.<br /> “Death was needed for the war to advance in the direction that it destined to go in.”<br /> <br /> Destined to go in?<br /> <br /> “My walking stick child!”'
$ws.Cells.Item(92, 27).Value = 'He leans upon the stick, He is old again, haggard and worn out. He walks slowly back to the rock where he sits for a moment before continuing on.<br /> <br /> “The Creator, or as we knew him: The Child, had learned something. His parents were alive. When things had started to get complicated, we intervened in the redirecting of the souls of his then deceased parents towards where you are now. You remember them. They were your parents.”<br /> <br /> He pauses and before you can really take in his words, and let them absorb, he continues:<br /> <br /> “The Child stated that the attack would be the best move, in order to spark the war. They had no idea that the target was not whom they assumed would be in the tent. This is why The Pirate Lords Son rages. This is why the world burns as it turns.”<br /> <br /> But its not the reason The Creator is who he is now, it just paints more of a story. One filled with suicide, pain and misery.<br /> <br /> It can twist and break a person.'
$ws.Cells.Item(92, 28).Value = 'The Smugglers Are Back!'

# Row 93
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = 'A Pawn in The Game'
$ws.Cells.Item(93, 3).Value = 'ThePoet'
$ws.Cells.Item(93, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(93, 5).Value = 'The Creators Choice'
$ws.Cells.Item(93, 6).Value = 1
$ws.Cells.Item(93, 8).Value = 4500
$ws.Cells.Item(93, 9).Value = 20000
$ws.Cells.Item(93, 11).Value = 10000
$ws.Cells.Item(93, 15).Value = 1000000000
$ws.Cells.Item(93, 16).Value = 3000
$ws.Cells.Item(93, 20).Value = 'Pirates and Merchants'
$ws.Cells.Item(93, 26).Value = 'The Poet stands beside you in a field. The one you crossed as a child to get to where you are now. The one that you played in. The House on the other side of the field lays in ruin. Burned and barley standing.<br /> <br /> “There are memories you have here?” The Poet asks. You say nothing.<br /> <br /> “That’s right, you are the child who’s parents are not your parents. The parents given to you. The ones …” Stolen.<br /> <br /> “False memories are nothing to cling to child.” His voice is soft and kind, but fills you with anger and hate.”'
$ws.Cells.Item(93, 27).Value = '“If you are done being sad about something that you never really had, the you should --” You cut him off and ask him where you even come from. If your parents are The Creators dead parents, then who are you?<br /> “A creation of celestial energy.”<br /> <br /> A tool to destroy The Creator?<br /> <br /> “Yes. A pawn in a game, where the riddles are the answers and you play because of your own morbid curiosity.”'
$ws.Cells.Item(93, 28).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(93, 29).Value = 'The Creators Choice'

# Row 94
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = 'Who Am I?'
$ws.Cells.Item(94, 3).Value = 'TheWitch'
$ws.Cells.Item(94, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(94, 7).Value = 'Lost Tome Of Spells'
$ws.Cells.Item(94, 8).Value = 45000
$ws.Cells.Item(94, 9).Value = 15000
$ws.Cells.Item(94, 11).Value = 15000
$ws.Cells.Item(94, 15).Value = 25000000
$ws.Cells.Item(94, 16).Value = 4000
$ws.Cells.Item(94, 20).Value = 'A Pawn in The Game'
$ws.Cells.Item(94, 26).Value = 'You stand alone before the hut where The Witch calls home. Alas she is no where to be seen. Her hut is empty, dark and cold. It feels warded with old magic.'
$ws.Cells.Item(94, 27).Value = 'She greets you with a smile as she sees you walking back up with a book she knows you killed for.<br /> <br /> “Child.” She states as she rises to greet you and take the book from you.<br /> <br /> You ask her, as you hand over the tome, who you really are. A pawn in a grand game of kill the corrupted and twisted soul?<br /> <br /> “No child. You are a free spirit of the earth. One that has transcended the depths of hell, the pits of purgatory. Have you seen into the belly of the beast yet child?”<br /> <br /> You look at her, evaluating her words and then shake your head.<br /> <br /> “Look deeper into the pit of the purgatories depths.”'
$ws.Cells.Item(94, 28).Value = 'The Smugglers Are Back!'

# Row 95
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = 'Dungeons Depths'
$ws.Cells.Item(95, 3).Value = 'DungeonMaster'
$ws.Cells.Item(95, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(95, 7).Value = 'Clean White Feather'
$ws.Cells.Item(95, 12).Value = 'Ratty Toy Bear'
$ws.Cells.Item(95, 15).Value = 50000000
$ws.Cells.Item(95, 16).Value = 4000
$ws.Cells.Item(95, 20).Value = 'Who Am I?'
$ws.Cells.Item(95, 26).Value = 'You descend into the madness of the purgatories dungeons seeking and item for the Dungeon Master. There are rumors of ring that glows with honesty and the change of “colors” what ever that means.<br /> <br /> You slaughter one beast after the other, searching for the single clean, white feather. Why a feather? Because it is soft?'
$ws.Cells.Item(95, 27).Value = 'You managed to find the Feather and bring it to the Dungeon Master, you seem confused as to why the odd request.<br /> <br /> “Because child, it was apart of this teddy bear, stuffed with feathers.”<br /> <br /> You hold the bear in your hand, ragged now, but the memories are there. Your mother gave it to you when you were younger. She made the bear her self and stuffed it with doves feathers. She made it with love. You <br /> remember that love. But what you remember is a lie.'
$ws.Cells.Item(95, 28).Value = 'The Smugglers Are Back!'

# Row 96
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = 'Hear of Home'
$ws.Cells.Item(96, 3).Value = 'ThePoet'
$ws.Cells.Item(96, 4).Value = 'The Smugglers Are Back!'
$ws.Cells.Item(96, 8).Value = 100000
$ws.Cells.Item(96, 9).Value = 100000
$ws.Cells.Item(96, 11).Value = 75000
$ws.Cells.Item(96, 16).Value = 5000
$ws.Cells.Item(96, 20).Value = 'Dungeons Depths'
$ws.Cells.Item(96, 24).Value = 2
$ws.Cells.Item(96, 26).Value = 'You stand in the ruins of your home and wait for the memories to come back.<br /> <br /> The Poet walks into the shambles of the burnt out house.<br /> <br /> You look at him and the back at what use to be the Kitchen, a place your mother would cook dinner, a place where your father would teach you how to gut an animal.<br /> <br /> A place where you ate with those you once loved. But you also remember the words. Being a pawn, being a tool to destroy a twisted and dark creature.<br /> <br /> It all fills you with rage and tears you can’t seem to hold back. Don’t want to.<br /> <br /> “Why the tears child?”<br /> <br /> You explain in angry tear laced words the hatred you have for him and others like him. For the manipulation, the lies and the pain.<br /> <br /> “You can be as angry as you please, you can throw a tantrum if it will help. But that won’t help the situation at hand. The Parents are locked away, behind a veil of memories.”<br /> <br /> You wont free them, you wont set them free. They are not your parents, they are The Creators. You explain you know this game.'
$ws.Cells.Item(96, 27).Value = '“Do you child?” The Poet Asks.<br /> <br /> “Do you really know this game?”<br /> <br /> He laughs.<br /> <br /> “The Child failed to ascend. His own parents failed to help him. He never wanted to find peace. He wanted to stay in a fantasy, stay in a world where he felt like he belonged. It corrupted people, it corrupted the ones he <br /> came to know. It was the downfall of The Red Hawks.”<br /> <br /> He stops and turns away from you. You absorb the information.<br /> <br /> “The Creator is a corruption. One that must be expunged.”'
$ws.Cells.Item(96, 28).Value = 'The Smugglers Are Back!'
